# feat: add 2022-Q4 data
#
# 1) "总计" (sheet1): insert a new first data row for 2022-Q4 and push the
#    existing quarters down by one row.
# 2) Insert a brand-new "2022-Q4" worksheet (positioned right after "总计",
#    before "2022-Q3") carrying the per-fund breakdown for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryData = @(
    @(0, "2022-Q4", 13, 3.68),
    @(1, "2022-Q3", 26, 5.37),
    @(2, "2022-Q2", 18, 4.92),
    @(3, "2022-Q1", 17, 2.91),
    @(4, "2021-Q4", 19, 8.57)
)

# Give the brand-new row 6 the same style as row 5 (bold/bordered index
# cell in column A) before the values are written, so it matches the
# formatting already used by every other row in this column.
$summary.Cells.Item(5, 1).Copy($summary.Cells.Item(6, 1))

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Borrow the bold/bordered header + index-column style already used by
# every other quarterly sheet (e.g. "2022-Q3") instead of re-deriving it.
$styleSource = $wb.Worksheets.Item("2022-Q3")
$styleSource.Range("B1:H1").Copy($q4.Range("B1:H1"))
$styleSource.Range("A2:A14").Copy($q4.Range("A2:A14"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$fundData = @(
    @(0, "002560", "诺安和鑫灵活配置混合", "31.46", "84.32", "2.93", "0.9218", 10),
    @(1, "010003", "景顺长城电子信息产业股票A", "18.37", "93.05", "4.54", "0.8340", 6),
    @(2, "001404", "招商移动互联网产业股票A", "13.83", "90.58", "4.77", "0.6597", 6),
    @(3, "015773", "招商移动互联网产业股票C", "8.53", "90.58", "4.77", "0.4069", 6),
    @(4, "010004", "景顺长城电子信息产业股票C", "5.90", "93.05", "4.54", "0.2679", 6),
    @(5, "008655", "招商科技创新混合A", "3.94", "93.96", "4.77", "0.1879", 7),
    @(6, "506001", "万家科创板 2 年定期开放混合", "6.18", "95.02", "2.98", "0.1842", 9),
    @(7, "008656", "招商科技创新混合C", "3.59", "93.96", "4.77", "0.1712", 7),
    @(8, "005844", "东方人工智能主题混合", "0.74", "94.64", "4.77", "0.0353", 10),
    @(9, "008300", "人保量化锐进混合A", "0.08", "90.63", "4.34", "0.0035", 4),
    @(10, "005629", "汇安趋势动力股票C", "0.08", "91.16", "3.94", "0.0032", 10),
    @(11, "008301", "人保量化锐进混合C", "0.04", "90.63", "4.34", "0.0017", 4),
    @(12, "005628", "汇安趋势动力股票A", "0.02", "91.16", "3.94", "0.0008", 10)
)

$r = 2
foreach ($row in $fundData) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    # Fund code, fund size, and position figures are stored as text
    # (mirrors the other quarterly sheets, e.g. "002560" keeps its
    # leading zero and "31.46" stays text instead of becoming a number).
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
